# Update "想去人数" (want-to-go count) figures in the F column on the
# "展览" (sheet 1) and "全部类型" (sheet 4) tabs, which mirror the same
# events. Values refreshed from a newer scrape (gh-pages data regen).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 (sheet 1)
$wsExhibit.Range("F3").Value  = 3250
$wsExhibit.Range("F4").Value  = 1681
$wsExhibit.Range("F5").Value  = 2325
$wsExhibit.Range("F16").Value = 8170
$wsExhibit.Range("F28").Value = 1070
$wsExhibit.Range("F30").Value = 1709
$wsExhibit.Range("F41").Value = 370

# 全部类型 (sheet 4) - same events, rows offset by +2
$wsAll.Range("F5").Value  = 3250
$wsAll.Range("F6").Value  = 1681
$wsAll.Range("F7").Value  = 2325
$wsAll.Range("F17").Value = 8170
$wsAll.Range("F30").Value = 1071
$wsAll.Range("F32").Value = 1709
$wsAll.Range("F42").Value = 370
